# Update "想去人数" (interest count) figures for several rows.
# The same underlying rows are duplicated on the "展览" sheet and on the
# "全部类型" aggregate sheet, so both need to be updated identically.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

# Map of row number -> new value for column F.
$updates = @{
    4  = 11042
    5  = 10227
    11 = 30
    13 = 9585
    17 = 10
    20 = 10857
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
